$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "本人收礼地址"
$ws.Range("A6").Value = "朋友完善地址"
$ws.Range("B5").Value = "address.html"
$ws.Range("B6").Value = "address2.html"

$ws.Range("E7").Select()
